$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 header row: rename "april 3 est" -> "estimate" (C15 "april 3 actuals" stays as-is)
$ws.Range("B15").Value = "estimate"

# New date header columns for April 4 and April 5, 2016
$ws.Range("D15").Value = 42464
$ws.Range("D15").NumberFormat = "d-mmm"

$ws.Range("E15").Value = 42465
$ws.Range("E15").NumberFormat = "mm-dd-yy"

# New actuals for "gallery" task row (row 18)
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = 9

# Move the active selection to match the new layout
$ws.Range("E20").Select()
